$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new time log entry in row 83 (previously blank), pushing the
# blank-row block down by one. The E column formula is a shared formula
# already defined for E4:E103, so setting the formula on E83 (and it
# already exists on E83 because the f-range covers it) - we just need to
# fill in A83:D83 and F83 values; E83's formula already lives in the
# sheet. We set it explicitly to be safe.

$ws.Range("A83").Value = 41924
$ws.Range("B83").Value = 0.82638888888888884
$ws.Range("C83").Value = 0.87638888888888899
$ws.Range("D83").Value = 20
$ws.Range("E83").Formula = "=IF(AND(NOT(ISBLANK(B83)),NOT(ISBLANK(C83))), (C83-B83) * 24 - D83/60, `"`")"
$ws.Range("F83").Value = "Coding"

# Excel materializes the shared formula's cached results for the
# previously-implicit blank rows (E84:E102) once the block is touched,
# the same way it already had for E103. Re-apply the shared formula
# explicitly so those cells get their own (empty-string) cached value.
for ($r = 84; $r -le 102; $r++) {
    $ws.Range("E$r").Formula = "=IF(AND(NOT(ISBLANK(B$r)),NOT(ISBLANK(C$r))), (C$r-B$r) * 24 - D$r/60, `"`")"
}

# Move the active selection down to D84 (mirrors the user tabbing to the
# next row after filling in row 83).
$ws.Range("D84").Select()

$wb.Save()
